$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in row 3 and row 4 for columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# and P (Precio $/Kg). The other columns are identical between the two
# rows, so swapping them would have no visible effect.

$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")
    $v3 = $cell3.Value2()
    $v4 = $cell4.Value2()
    $cell3.Value = $v4
    $cell4.Value = $v3
}

$wb.Save()
